$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$words = @("버스","노선","퇴근","선","출근","시간","없","분","이용","정류장","운행","많","동","부탁","도착","시","출발","사람","통근","역","요청","인원","검토","경우","변경","건의","생각","차","출퇴근","감사","앞","동백","신설","증차","추가","문의","정차","영통","답변","회사","관련","안녕","탑승","수원","정도","불편","셔틀","차량","사항","잠실","청주","용인","위치","좋","현재","글","구리","이천","자리","아침","운영","중","지역","하차","배차","안양","인천","전","쪽","신","직원","아파트","수고","금요일","대","조사","양재","후","분당","확인","조정","마을","평일","공지","주말","거리","교대","때문","행","하세","안녕하","문제","이상","곳","일","기사님","담당자","상황","천호","태")
$counts = @(4767,2424,1903,1763,1533,1504,1157,1079,1072,965,960,927,809,746,725,721,717,701,701,682,661,652,598,579,571,565,559,519,501,497,493,483,475,471,457,457,454,454,453,452,445,444,444,424,420,416,409,394,390,378,378,363,354,353,348,327,313,309,303,303,302,296,295,293,285,269,263,261,258,256,256,254,253,249,244,242,242,240,240,239,234,227,226,220,214,212,212,208,204,203,202,202,201,201,201,200,199,198,194,193)

for ($i = 0; $i -lt $words.Length; $i++) {
    $r = 2 + $i
    $ws.Cells.Item($r, 2).Value = $words[$i]
    $ws.Cells.Item($r, 3).Value = $counts[$i]
}
